$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the whole data range is treated as text so numeric-looking
# strings (e.g. "1.001", "29.178.44") are NOT auto-converted to numbers,
# matching the source data which stores these as plain text.
$dataRange = $ws.Range('B2:E51')
$dataRange.NumberFormat = '@'

# Row 2
$ws.Range('D2').Value = '29.178.44'
$ws.Range('E2').Value = '  +0.42%  '

# Row 3
$ws.Range('D3').Value = '1.841.53'
$ws.Range('E3').Value = '  +0.39%  '

# Row 4
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  +0.12%  '

# Row 5
$ws.Range('D5').Value = '244.18'
$ws.Range('E5').Value = '  -0.29%  '

# Row 6
$ws.Range('D6').Value = '0.6257'
$ws.Range('E6').Value = '  -1.39%  '

# Row 7
$ws.Range('E7').Value = '  +0.09%  '

# Row 8
$ws.Range('D8').Value = '0.07546'
$ws.Range('E8').Value = '  -0.60%  '

# Row 9
$ws.Range('D9').Value = '0.2951'
$ws.Range('E9').Value = '  -0.04%  '

# Row 10
$ws.Range('D10').Value = '23.37'
$ws.Range('E10').Value = '  +2.33%  '

# Row 11
$ws.Range('D11').Value = '0.07712'
$ws.Range('E11').Value = '  -0.50%  '

# Row 12
$ws.Range('D12').Value = '1.853.23'
$ws.Range('E12').Value = '  +1.38%  '

# Row 13
$ws.Range('D13').Value = '5.032'
$ws.Range('E13').Value = '  +0.61%  '

# Row 14
$ws.Range('D14').Value = '0.6795'
$ws.Range('E14').Value = '  +1.15%  '

# Row 15
$ws.Range('D15').Value = '83.28'
$ws.Range('E15').Value = '  -0.01%  '

# Row 16
$ws.Range('D16').Value = '0.000009316'
$ws.Range('E16').Value = '  -5.39%  '

# Row 17
$ws.Range('D17').Value = '5.992'
$ws.Range('E17').Value = '  -2.05%  '

# Row 18
$ws.Range('D18').Value = '29.174.27'
$ws.Range('E18').Value = '  +0.38%  '

# Row 19
$ws.Range('D19').Value = '2.086.92'
$ws.Range('E19').Value = '  +0.08%  '

# Row 20
$ws.Range('D20').Value = '232.47'
$ws.Range('E20').Value = '  +2.42%  '

# Row 21
$ws.Range('D21').Value = '12.73'
$ws.Range('E21').Value = '  +1.24%  '

# Row 22
$ws.Range('E22').Value = '  +0.22%  '

# Row 23
$ws.Range('D23').Value = '7.190'
$ws.Range('E23').Value = '  -0.49%  '

# Row 24
$ws.Range('E24').Value = '  +0.10%  '

# Row 25
$ws.Range('D25').Value = '160.82'
$ws.Range('E25').Value = '  +0.19%  '

# Row 26
$ws.Range('D26').Value = '0.1407'
$ws.Range('E26').Value = '  +0.18%  '

# Row 27
$ws.Range('D27').Value = '8.578'
$ws.Range('E27').Value = '  +0.33%  '

# Row 28
$ws.Range('E28').Value = '  +0.08%  '

# Row 29
$ws.Range('E29').Value = '  -0.32%  '

# Row 30
$ws.Range('D30').Value = '4.198'
$ws.Range('E30').Value = '  +1.69%  '

# Row 31
$ws.Range('D31').Value = '4.157'
$ws.Range('E31').Value = '  +2.43%  '

# Row 32
$ws.Range('D32').Value = '0.05575'
$ws.Range('E32').Value = '  +3.41%  '

# Row 33
$ws.Range('D33').Value = '1.208'
$ws.Range('E33').Value = '  +0.15%  '

# Row 34
$ws.Range('D34').Value = '0.7503'
$ws.Range('E34').Value = '  +0.22%  '

# Row 35
$ws.Range('E35').Value = '  -0.47%  '

# Row 36
$ws.Range('D36').Value = '1.150'
$ws.Range('E36').Value = '  +0.59%  '

# Row 37
$ws.Range('E37').Value = '  +0.00%  '

# Row 38
$ws.Range('D38').Value = '1.240.33'
$ws.Range('E38').Value = '  +0.35%  '

# Row 39
$ws.Range('D39').Value = '2.775'
$ws.Range('E39').Value = '  +0.55%  '

# Row 40
$ws.Range('D40').Value = '0.01796'
$ws.Range('E40').Value = '  -0.07%  '

# Row 41
$ws.Range('D41').Value = '6.626'
$ws.Range('E41').Value = '  -0.05%  '

# Row 42
$ws.Range('D42').Value = '0.9030'
$ws.Range('E42').Value = '  -0.05%  '

# Row 43
$ws.Range('E43').Value = '  -0.01%  '

# Row 44
$ws.Range('D44').Value = '102.40'
$ws.Range('E44').Value = '  -0.21%  '

# Row 45
$ws.Range('D45').Value = '67.00'
$ws.Range('E45').Value = '  +3.16%  '

# Row 46
$ws.Range('D46').Value = '1.986.89'
$ws.Range('E46').Value = '  +0.11%  '

# Row 47
$ws.Range('E47').Value = '  -1.25%  '

# Row 48
$ws.Range('D48').Value = '0.5097'
$ws.Range('E48').Value = '  -0.30%  '

# Row 49
$ws.Range('D49').Value = '0.4099'
$ws.Range('E49').Value = '  +0.03%  '

# Row 50
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '9.113'
$ws.Range('E50').Value = '  +0.53%  '

# Row 51
$ws.Range('B51').Value = 'XinFinNetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$ws.Range('D51').Value = '0.07357'
$ws.Range('E51').Value = '  +16.18%  '

# Restore the default "Normal" style on the data range so no stray
# number-format/style indices are left behind on the written cells.
$dataRange.Style = 'Normal'

Write-Output "Applied all cell updates"